$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("morpho/morpho-v1-main/src/aave-v2/EntryPositionsManager.sol", 296, 1, 2, 0),
    @("morpho/morpho-v1-main/src/aave-v2/ExitPositionsManager.sol", 715, 1, 2, 0),
    @("morpho/morpho-v1-main/src/aave-v2/IncentivesVault.sol", 151, 1, 3, 3),
    @("morpho/morpho-v1-main/src/aave-v2/InterestRatesManager.sol", 189, 1, 3, 3),
    @("morpho/morpho-v1-main/src/aave-v2/MatchingEngine.sol", 354, 1, 1, 0),
    @("morpho/morpho-v1-main/src/aave-v2/Morpho.sol", 245, 1, 1, 0),
    @("morpho/morpho-v1-main/src/aave-v2/MorphoGovernance.sol", 531, 1, 1, 1),
    @("morpho/morpho-v1-main/src/aave-v2/MorphoStorage.sol", 87, 1, 7, 3),
    @("morpho/morpho-v1-main/src/aave-v2/MorphoUtils.sol", 338, 1, 5, 4),
    @("morpho/morpho-v1-main/src/aave-v2/PositionsManagerUtils.sol", 100, 1, 2, 1),
    @("morpho/morpho-v1-main/src/aave-v2/interfaces/IEntryPositionsManager.sol", 19, 1, 0, 0),
    @("morpho/morpho-v1-main/src/aave-v2/interfaces/IExitPositionsManager.sol", 30, 1, 0, 0),
    @("morpho/morpho-v1-main/src/aave-v2/interfaces/IGetterUnderlyingAsset.sol", 7, 1, 0, 0),
    @("morpho/morpho-v1-main/src/aave-v2/interfaces/IIncentivesVault.sol", 29, 1, 1, 0),
    @("morpho/morpho-v1-main/src/aave-v2/interfaces/IInterestRatesManager.sol", 7, 1, 0, 0),
    @("morpho/morpho-v1-main/src/aave-v2/interfaces/IMorpho.sol", 103, 1, 8, 0),
    @("morpho/morpho-v1-main/src/aave-v2/interfaces/IOracle.sol", 7, 1, 0, 0),
    @("morpho/morpho-v1-main/src/aave-v2/interfaces/IRewardsManager.sol", 27, 1, 1, 0),
    @("morpho/morpho-v1-main/src/aave-v2/interfaces/aave/IAToken.sol", 91, 1, 2, 0),
    @("morpho/morpho-v1-main/src/aave-v2/interfaces/aave/IAaveIncentivesController.sol", 154, 2, 0, 0),
    @("morpho/morpho-v1-main/src/aave-v2/interfaces/aave/IERC20.sol", 81, 1, 0, 0),
    @("morpho/morpho-v1-main/src/aave-v2/interfaces/aave/ILendingPool.sol", 418, 1, 2, 0),
    @("morpho/morpho-v1-main/src/aave-v2/interfaces/aave/ILendingPoolAddressesProvider.sol", 61, 1, 0, 0),
    @("morpho/morpho-v1-main/src/aave-v2/interfaces/aave/IPriceOracleGetter.sol", 13, 1, 0, 0),
    @("morpho/morpho-v1-main/src/aave-v2/interfaces/aave/IScaledBalanceToken.sol", 27, 1, 0, 0),
    @("morpho/morpho-v1-main/src/aave-v2/interfaces/aave/IVariableDebtToken.sol", 79, 1, 1, 0),
    @("morpho/morpho-v1-main/src/aave-v2/interfaces/lido/ILido.sol", 9, 1, 0, 0),
    @("morpho/morpho-v1-main/src/aave-v2/lens/IndexesLens.sol", 126, 1, 3, 0),
    @("morpho/morpho-v1-main/src/aave-v2/lens/Lens.sol", 108, 1, 1, 0),
    @("morpho/morpho-v1-main/src/aave-v2/lens/LensStorage.sol", 50, 1, 8, 3),
    @("morpho/morpho-v1-main/src/aave-v2/lens/MarketsLens.sol", 235, 1, 1, 0),
    @("morpho/morpho-v1-main/src/aave-v2/lens/RatesLens.sol", 469, 1, 2, 0),
    @("morpho/morpho-v1-main/src/aave-v2/lens/UsersLens.sol", 483, 1, 1, 1),
    @("morpho/morpho-v1-main/src/aave-v2/lens/interfaces/ILens.sol", 257, 1, 3, 0),
    @("morpho/morpho-v1-main/src/aave-v2/libraries/InterestRatesModel.sol", 207, 1, 1, 3),
    @("morpho/morpho-v1-main/src/aave-v2/libraries/Types.sol", 105, 1, 0, 0),
    @("morpho/morpho-v1-main/src/aave-v2/libraries/aave/DataTypes.sol", 54, 1, 0, 0),
    @("morpho/morpho-v1-main/src/aave-v2/libraries/aave/Errors.sol", 120, 1, 0, 0),
    @("morpho/morpho-v1-main/src/aave-v2/libraries/aave/ReserveConfiguration.sol", 376, 1, 2, 0),
    @("morpho/morpho-v1-main/src/aave-v2/libraries/aave/UserConfiguration.sol", 17, 1, 1, 0),
    @("morpho/morpho-v1-main/src/common/rewards-distribution/RewardsDistributor.sol", 104, 1, 0, 3),
    @("morpho/morpho-v1-main/src/common/test/FakeToken.sol", 13, 1, 0, 1),
    @("morpho/morpho-v1-main/src/compound/IncentivesVault.sol", 156, 1, 3, 3),
    @("morpho/morpho-v1-main/src/compound/InterestRatesManager.sol", 171, 1, 3, 1),
    @("morpho/morpho-v1-main/src/compound/MatchingEngine.sol", 386, 1, 1, 0),
    @("morpho/morpho-v1-main/src/compound/Morpho.sol", 247, 1, 1, 0),
    @("morpho/morpho-v1-main/src/compound/MorphoGovernance.sol", 487, 1, 1, 0),
    @("morpho/morpho-v1-main/src/compound/MorphoStorage.sol", 76, 1, 6, 3),
    @("morpho/morpho-v1-main/src/compound/MorphoUtils.sol", 227, 1, 2, 3),
    @("morpho/morpho-v1-main/src/compound/PositionsManager.sol", 1032, 1, 3, 0),
    @("morpho/morpho-v1-main/src/compound/RewardsManager.sol", 262, 1, 3, 1),
    @("morpho/morpho-v1-main/src/compound/interfaces/IIncentivesVault.sol", 29, 1, 1, 0),
    @("morpho/morpho-v1-main/src/compound/interfaces/IInterestRatesManager.sol", 7, 1, 0, 0),
    @("morpho/morpho-v1-main/src/compound/interfaces/IMorpho.sol", 88, 1, 5, 0),
    @("morpho/morpho-v1-main/src/compound/interfaces/IOracle.sol", 7, 1, 0, 0),
    @("morpho/morpho-v1-main/src/compound/interfaces/IPositionsManager.sol", 44, 1, 0, 0),
    @("morpho/morpho-v1-main/src/compound/interfaces/IRewardsManager.sol", 39, 1, 1, 0),
    @("morpho/morpho-v1-main/src/compound/interfaces/IWETH.sol", 9, 1, 0, 0),
    @("morpho/morpho-v1-main/src/compound/interfaces/compound/ICompound.sol", 380, 6, 0, 0),
    @("morpho/morpho-v1-main/src/compound/lens/IndexesLens.sol", 159, 1, 2, 0),
    @("morpho/morpho-v1-main/src/compound/lens/Lens.sol", 97, 1, 1, 0),
    @("morpho/morpho-v1-main/src/compound/lens/LensStorage.sol", 48, 1, 4, 5),
    @("morpho/morpho-v1-main/src/compound/lens/MarketsLens.sol", 180, 1, 1, 0),
    @("morpho/morpho-v1-main/src/compound/lens/RatesLens.sol", 453, 1, 1, 0),
    @("morpho/morpho-v1-main/src/compound/lens/RewardsLens.sol", 181, 1, 1, 0),
    @("morpho/morpho-v1-main/src/compound/lens/UsersLens.sol", 471, 1, 1, 0),
    @("morpho/morpho-v1-main/src/compound/lens/interfaces/ILens.sol", 291, 1, 3, 0),
    @("morpho/morpho-v1-main/src/compound/libraries/CompoundMath.sol", 53, 1, 0, 0),
    @("morpho/morpho-v1-main/src/compound/libraries/InterestRatesModel.sol", 203, 1, 1, 3),
    @("morpho/morpho-v1-main/src/compound/libraries/Types.sol", 94, 1, 0, 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

Write-Output "Updated $($data.Length) rows"